$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "Structure N"/"TSx,y" headers to the new Delta-G-RRS labels ---
$ws.Range("B1").Value = "ΔGRRS(2)"
$ws.Range("C1").Value = "ΔGRRS(3)"
$ws.Range("D1").Value = "ΔGRRS(TS1)"
$ws.Range("E1").Value = "ΔGRRS(4)"
$ws.Range("F1").Value = "ΔGRRS(5)"
$ws.Range("G1").Value = "ΔGRRS(TS2) "
$ws.Range("H1").Value = "ΔGRRS(6)"
$ws.Range("I1").Value = "ΔGRRS(TS3)  "
$ws.Range("J1").Value = "ΔGRRS(7)"
$ws.Range("K1").Value = "ΔGRRS(TS4)  "

# --- 2. B1 drops its special Cantarell/bordered/0.00 formatting and becomes a
#        plain, centered, General-format cell (matching the rest of the header row) ---
$ws.Range("B1").Font.Name = "Arial"
$ws.Range("B1").Font.Size = 10
$ws.Range("B1").Borders.LineStyle = -4142
$ws.Range("B1").NumberFormat = "General"

# --- 3. C1,D1,E1,F1,H1,J1 switch from left-aligned "0.00" to centered "General"
#        (G1,I1,K1 stay as they were: left-aligned, "0.00") ---
$centeredHeaders = @("C1","D1","E1","F1","H1","J1")
foreach ($addr in $centeredHeaders) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).HorizontalAlignment = -4108
}

# --- 4. Row 1 height tweak ---
$ws.Rows.Item(1).RowHeight = 13.55

# --- 5. Selection/view bookkeeping ---
$ws.Range("P16").Select()
